# B1--and-B2-PowerPoint.pptx edit
#
# The authored change (per the canonical-OOXML diff) updates the table-style
# applied to the financial-documents table on slide 5 (the slide with the
# "B1- TYPES OF FINANCIAL DOCUMENTS" title) from the default "Table_0" style
#   {B45DC1D0-0F84-4E7F-B27C-459705C21386}
# to the built-in PowerPoint table style
#   {B6147149-64DA-4463-A9C5-0B8D8A7B7426}
#
# Table styles can't be changed by assigning Table.Style (that property is
# read-only in the PowerPoint object model - PowerPoint raises "Table styles
# cannot be assigned through a property ... call Table.ApplyStyle(...)
# instead"); Table.ApplyStyle(styleId) is the supported call, so we look up
# the table shape on every slide and (re)apply the new style GUID wherever a
# table is found using the previous style.

$p = $ppt.ActivePresentation

$oldStyleId = "{B45DC1D0-0F84-4E7F-B27C-459705C21386}"
$newStyleId = "{B6147149-64DA-4463-A9C5-0B8D8A7B7426}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}
